# Append a new logbook entry as row 26, reflecting an additional
# processing run recorded after updating the input source option.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

$ws.Cells.Item($row, 1).Value  = 25
$ws.Cells.Item($row, 2).Value  = "2025-11-29 16:40:54"
$ws.Cells.Item($row, 3).Value  = "A873-150925-CHK-Y06"
$ws.Cells.Item($row, 4).Value  = 891
$ws.Cells.Item($row, 5).Value  = 891
$ws.Cells.Item($row, 6).Value  = 881
$ws.Cells.Item($row, 7).Value  = 0
$ws.Cells.Item($row, 8).Value  = 9
$ws.Cells.Item($row, 9).Value  = 0
$ws.Cells.Item($row, 10).Value = 1
$ws.Cells.Item($row, 11).Value = 631
$ws.Cells.Item($row, 12).Value = $false
$ws.Cells.Item($row, 13).Value = 10
$ws.Cells.Item($row, 14).Value = 1.12
$ws.Cells.Item($row, 15).Value = 2.2
